# Apply the "Update countries & provincias Spain" edit:
#  - A1: refresh the "Datos actualizados..." timestamp string
#  - A4:H215: re-sorted country table (by "Casos totales" desc) with updated counts
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'Datos actualizados a 10 de Abril de 2020 a las 00:22'

$rows = @(
  @('Estados Unidos', 465088, 30161, 25139, 423439, 9823, 1722, 16510),
  @('España', 153222, 5002, 52165, 85610, 7371, 655, 15447),
  @('Italia', 143626, 4204, 28470, 96877, 3605, 610, 18279),
  @('Francia', 117749, 4799, 23206, 82333, 7066, 1341, 12210),
  @('Alemania', 117712, 4416, 52407, 62776, 4895, 180, 2529),
  @('China', 81865, 0, 77370, 1160, 176, 0, 3335),
  @('Iran', 66220, 1634, 32309, 29801, 3918, 117, 4110),
  @('Reino Unido', 65077, 4344, 135, 56964, 1559, 881, 7978),
  @('Turquia', 42282, 4056, 2142, 39232, 1552, 96, 908),
  @('Belgica', 24983, 1580, 5164, 17296, 1285, 283, 2523),
  @('Suiza', 24051, 771, 10600, 12503, 386, 53, 948),
  @('Paises Bajos', 21762, 1213, 250, 19116, 1424, 148, 2396),
  @('Canada', 20703, 1265, 5218, 14981, 518, 77, 504),
  @('Brasil', 18092, 1904, 173, 16969, 296, 130, 950),
  @('Portugal', 13956, 815, 205, 13342, 241, 29, 409),
  @('Austria', 13244, 302, 5240, 7709, 266, 22, 295),
  @('Corea del Sur', 10423, 39, 6973, 3246, 55, 4, 204),
  @('Rusia', 10131, 1459, 698, 9357, 8, 13, 76),
  @('Israel', 9968, 564, 1011, 8871, 166, 13, 86),
  @('Suecia', 9141, 722, 205, 8143, 719, 106, 793),
  @('India', 6725, 809, 620, 5879, 0, 48, 226),
  @('Irlanda', 6574, 500, 25, 6286, 165, 28, 263),
  @('Noruega', 6162, 120, 32, 6022, 78, 7, 108),
  @('Australia', 6104, 52, 2987, 3066, 81, 1, 51),
  @('Chile', 5972, 426, 1274, 4641, 360, 9, 57),
  @('Dinamarca', 5635, 233, 1736, 3662, 120, 19, 237),
  @('Polonia', 5575, 370, 284, 5117, 160, 15, 174),
  @('Chequia', 5467, 155, 301, 5054, 96, 13, 112),
  @('Peru', 5256, 914, 1438, 3680, 124, 17, 138),
  @('Rumania', 5202, 441, 647, 4307, 178, 28, 248),
  @('Ecuador', 4965, 515, 339, 4354, 139, 30, 272),
  @('Japon', 4667, 0, 632, 3941, 99, 0, 94),
  @('Pakistan', 4489, 226, 572, 3852, 31, 4, 65),
  @('Malasia', 4228, 109, 1608, 2553, 72, 2, 67),
  @('Filipinas', 4076, 206, 124, 3749, 1, 21, 203),
  @('Indonesia', 3293, 337, 252, 2761, 0, 40, 280),
  @('Arabia Saudita', 3287, 355, 666, 2577, 41, 3, 44),
  @('Mexico', 3181, 396, 633, 2374, 89, 33, 174),
  @('Luxemburgo', 3115, 81, 500, 2563, 30, 6, 52),
  @('Emiratos Arabes Unidos', 2990, 331, 268, 2708, 1, 2, 14),
  @('Serbia', 2867, 201, 118, 2683, 127, 1, 66),
  @('Finlandia', 2605, 118, 300, 2263, 82, 2, 42),
  @('Panama', 2528, 0, 16, 2449, 101, 0, 63),
  @('Tailandia', 2423, 54, 940, 1451, 61, 2, 32),
  @('Catar', 2376, 166, 206, 2164, 37, 0, 6),
  @('Republica Dominicana', 2349, 238, 80, 2151, 147, 10, 118),
  @('Colombia', 2223, 169, 174, 1980, 85, 14, 69),
  @('Grecia', 1955, 71, 269, 1599, 79, 4, 87),
  @('Sudafrica', 1934, 89, 95, 1821, 7, 0, 18),
  @('Singapur', 1910, 287, 460, 1444, 29, 0, 6),
  @('Ucrania', 1892, 224, 45, 1790, 33, 5, 57),
  @('Argentina', 1795, 0, 365, 1359, 96, 6, 71),
  @('Egipto', 1699, 139, 348, 1233, 0, 15, 118),
  @('Argelia', 1666, 94, 347, 1084, 46, 30, 235),
  @('Islandia', 1648, 32, 688, 954, 11, 0, 6),
  @('Bielorrusia', 1486, 420, 139, 1331, 55, 3, 16),
  @('Croacia', 1407, 64, 219, 1168, 34, 1, 20),
  @('Marruecos', 1374, 99, 109, 1168, 1, 4, 97),
  @('Moldavia', 1289, 115, 50, 1210, 80, 2, 29),
  @('Nueva Zelanda', 1239, 29, 317, 921, 4, 0, 1),
  @('Irak', 1232, 30, 496, 667, 0, 0, 69),
  @('Estonia', 1207, 22, 83, 1100, 9, 0, 24),
  @('Eslovenia', 1124, 33, 128, 953, 34, 3, 43),
  @('Hungria', 980, 85, 96, 818, 17, 8, 66),
  @('Hong Kong', 974, 13, 293, 677, 14, 0, 4),
  @('Lituania', 955, 43, 8, 931, 21, 1, 16),
  @('Azerbaiyan', 926, 104, 101, 816, 27, 1, 9),
  @('Armenia', 921, 40, 138, 773, 30, 1, 10),
  @('Kuwait', 910, 55, 111, 798, 22, 0, 1),
  @('Barein', 887, 64, 519, 363, 3, 0, 5),
  @('Bosnia y Herzegovina', 858, 54, 101, 722, 4, 1, 35),
  @('Kazajistan', 781, 54, 60, 713, 21, 1, 8),
  @('Camerun', 730, 0, 60, 660, 0, 0, 10),
  @('Crucero', 712, 0, 619, 82, 10, 0, 11),
  @('Eslovaquia', 701, 19, 23, 676, 5, 0, 2),
  @('Republica de Macedonia', 663, 46, 37, 596, 15, 1, 30),
  @('Tunez', 643, 15, 25, 593, 78, 1, 25),
  @('Bulgaria', 618, 25, 48, 546, 32, 0, 24),
  @('Letonia', 589, 12, 16, 570, 3, 1, 3),
  @('Principado de Andorra', 583, 19, 58, 500, 17, 2, 25),
  @('Uzbekistan', 582, 37, 38, 541, 8, 0, 3),
  @('Libano', 582, 6, 67, 496, 29, 0, 19),
  @('Republica de Chipre', 564, 38, 53, 501, 13, 1, 10),
  @('Costa Rica', 539, 37, 30, 506, 13, 0, 3),
  @('Cuba', 515, 58, 28, 472, 15, 3, 15),
  @('Afganistan', 484, 40, 32, 437, 0, 1, 15),
  @('Oman', 457, 38, 109, 345, 3, 1, 3),
  @('Uruguay', 456, 0, 192, 257, 14, 0, 7),
  @('Costa de Marfil', 444, 60, 52, 389, 0, 0, 3),
  @('Burkina Faso', 443, 29, 146, 273, 0, 1, 24),
  @('Niger', 410, 68, 40, 359, 0, 0, 11),
  @('Albania', 409, 9, 165, 221, 7, 1, 23),
  @('Taiwan', 380, 1, 80, 295, 0, 0, 5),
  @('Jordania', 372, 14, 161, 204, 5, 1, 7),
  @('Reunion', 362, 0, 40, 322, 4, 0, 0),
  @('Honduras', 343, 31, 6, 314, 10, 1, 23),
  @('Malta', 337, 38, 16, 319, 4, 1, 2),
  @('San Marino', 333, 25, 49, 250, 14, 0, 34),
  @('Banglades', 330, 112, 33, 276, 1, 1, 21),
  @('Mauricio', 314, 41, 23, 284, 3, 0, 7),
  @('Ghana', 313, 0, 34, 273, 2, 0, 6),
  @('Nigeria', 288, 12, 51, 230, 2, 1, 7),
  @('Kirguistan', 280, 10, 35, 241, 5, 0, 4),
  @('Bolivia', 264, 54, 2, 244, 3, 3, 18),
  @('Estado de Palestina', 263, 0, 44, 218, 0, 0, 1),
  @('Vietnam', 255, 4, 128, 127, 8, 0, 0),
  @('Montenegro', 252, 4, 4, 246, 7, 0, 2),
  @('Senegal', 250, 6, 123, 125, 1, 0, 2),
  @('Georgia', 218, 7, 51, 164, 6, 0, 3),
  @('Guinea', 194, 30, 11, 183, 0, 0, 0),
  @('Sri Lanka', 190, 1, 49, 134, 5, 0, 7),
  @('Isla de Man', 190, 32, 92, 97, 12, 0, 1),
  @('Kenia', 184, 5, 12, 165, 2, 1, 7),
  @('Mayotte', 184, 0, 22, 160, 3, 0, 2),
  @('Islas Feroe', 184, 0, 136, 48, 0, 0, 0),
  @('Consejo Danes para los Refugiados', 180, 0, 9, 153, 0, 0, 18),
  @('Venezuela', 171, 4, 84, 78, 6, 0, 9),
  @('Martinica', 154, 0, 50, 98, 19, 0, 6),
  @('Guadalupe', 141, 0, 43, 90, 13, 0, 8),
  @('Republica de Yibuti', 135, 0, 25, 110, 0, 0, 0),
  @('Brunei', 135, 0, 92, 42, 3, 0, 1),
  @('Paraguay', 124, 5, 18, 101, 1, 0, 5),
  @('Gibraltar', 123, 3, 60, 63, 1, 0, 0),
  @('Camboya', 119, 2, 62, 57, 1, 0, 0),
  @('Ruanda', 110, 0, 7, 103, 0, 0, 0),
  @('Trinidad yTobago', 109, 2, 1, 100, 0, 0, 8),
  @('El Salvador', 103, 10, 14, 83, 3, 1, 6),
  @('Guatemala', 95, 8, 17, 75, 3, 0, 3),
  @('Madagascar', 93, 0, 11, 82, 1, 0, 0),
  @('Monaco', 84, 3, 5, 78, 4, 0, 1),
  @('Guayana Francesa', 83, 0, 43, 40, 1, 0, 0),
  @('Aruba', 82, 5, 20, 62, 0, 0, 0),
  @('Liechtenstein', 78, 0, 55, 22, 0, 0, 1),
  @('Mali', 74, 15, 22, 45, 0, 0, 7),
  @('Togo', 73, 3, 24, 46, 0, 0, 3),
  @('Barbados', 66, 3, 11, 52, 4, 0, 3),
  @('Jamaica', 63, 0, 12, 47, 0, 0, 4),
  @('Congo', 60, 15, 5, 50, 0, 0, 5),
  @('Etiopia', 56, 1, 4, 50, 2, 0, 2),
  @('Uganda', 53, 0, 0, 53, 0, 0, 0),
  @('Polinesia Francesa', 51, 0, 0, 51, 0, 0, 0),
  @('Bermudas', 48, 9, 25, 19, 0, 1, 4),
  @('Islas Caimanes', 45, 0, 6, 38, 0, 0, 1),
  @('Macao', 45, 0, 10, 35, 1, 0, 0),
  @('Gabon', 44, 10, 1, 42, 0, 0, 1),
  @('San Martin (Parte Holandesa)', 43, 3, 1, 36, 2, 0, 6),
  @('Bahamas', 40, 0, 5, 28, 1, 0, 7),
  @('Puerto Rico', 39, 0, 1, 36, 0, 0, 2),
  @('Zambia', 39, 0, 24, 14, 1, 0, 1),
  @('Guyana', 37, 0, 8, 23, 4, 0, 6),
  @('Guinea-Bisau', 36, 3, 0, 36, 0, 0, 0),
  @('Eritrea', 33, 0, 0, 33, 0, 0, 0),
  @('Guam', 32, 0, 0, 31, 0, 0, 1),
  @('San Martin (Parte Francesa)', 32, 0, 7, 23, 6, 0, 2),
  @('Liberia', 31, 0, 3, 24, 0, 0, 4),
  @('Haiti', 30, 3, 0, 28, 0, 1, 2),
  @('Benin', 26, 0, 5, 20, 0, 0, 1),
  @('Tanzania', 25, 0, 5, 19, 0, 0, 1),
  @('Libia', 24, 3, 8, 15, 0, 0, 1),
  @('Birmania', 23, 1, 2, 18, 0, 0, 3),
  @('Antigua y Barbuda', 19, 0, 0, 17, 1, 0, 2),
  @('Angola', 19, 0, 2, 15, 0, 0, 2),
  @('Siria', 19, 0, 4, 13, 0, 0, 2),
  @('Maldivas', 19, 0, 13, 6, 0, 0, 0),
  @('Nueva Caledonia', 18, 0, 1, 17, 0, 0, 0),
  @('Guinea Ecuatorial', 18, 0, 3, 15, 0, 0, 0),
  @('Islas Virgenes de los Estados Unidos', 17, 0, 0, 17, 0, 0, 0),
  @('Mozambique', 17, 0, 1, 16, 0, 0, 0),
  @('Laos', 16, 1, 0, 16, 0, 0, 0),
  @('Namibia', 16, 0, 3, 13, 0, 0, 0),
  @('Mongolia', 16, 0, 4, 12, 0, 0, 0),
  @('Fiyi', 15, 0, 0, 15, 0, 0, 0),
  @('Dominica', 15, 0, 1, 14, 0, 0, 0),
  @('Sudan', 15, 1, 2, 11, 0, 0, 2),
  @('Santa Lucia', 14, 0, 1, 13, 0, 0, 0),
  @('Curazao', 14, 0, 7, 6, 0, 0, 1),
  @('Botsuana', 13, 7, 0, 12, 0, 0, 1),
  @('Granada', 12, 0, 0, 12, 2, 0, 0),
  @('San Vicente y las Granadinas', 12, 4, 1, 11, 0, 0, 0),
  @('Somalia', 12, 0, 1, 10, 0, 0, 1),
  @('Suazilandia', 12, 0, 7, 5, 0, 0, 0),
  @('San Cristobal y Nieves', 11, 0, 0, 11, 0, 0, 0),
  @('Seychelles', 11, 0, 0, 11, 0, 0, 0),
  @('Republica del Chad', 11, 1, 2, 9, 0, 0, 0),
  @('Zimbabue', 11, 0, 0, 8, 0, 0, 3),
  @('Groenlandia', 11, 0, 11, 0, 0, 0, 0),
  @('Surinam', 10, 0, 4, 5, 0, 0, 1),
  @('Belice', 9, 1, 0, 8, 1, 0, 1),
  @('Nepal', 9, 0, 1, 8, 0, 0, 0),
  @('Montserrat', 9, 0, 0, 7, 0, 0, 2),
  @('Republica de Africa Central', 8, 0, 0, 8, 0, 0, 0),
  @('Islas Turcas y Caicos', 8, 0, 0, 7, 0, 0, 1),
  @('Malaui', 8, 0, 0, 7, 1, 0, 1),
  @('Santa Sede', 8, 0, 2, 6, 0, 0, 0),
  @('Sierra Leona', 7, 0, 0, 7, 0, 0, 0),
  @('Nicaragua', 7, 1, 0, 6, 0, 0, 1),
  @('Cabo Verde', 7, 0, 1, 5, 0, 0, 1),
  @('Mauritania', 7, 1, 2, 4, 0, 0, 1),
  @('San Bartolome', 6, 0, 1, 5, 0, 0, 0),
  @('Islas Malvinas', 5, 0, 1, 4, 0, 0, 0),
  @('Butan', 5, 0, 2, 3, 0, 0, 0),
  @('Santo Tome y Principe', 4, 0, 0, 4, 0, 0, 0),
  @('Sahara Occidental', 4, 0, 0, 4, 0, 0, 0),
  @('Gambia', 4, 0, 2, 1, 0, 0, 1),
  @('Anguila', 3, 0, 0, 3, 0, 0, 0),
  @('Islas Virgenes Britanicas', 3, 0, 0, 3, 0, 0, 0),
  @('Sudan del Sur', 3, 1, 0, 3, 0, 0, 0),
  @('Burundi', 3, 0, 0, 3, 0, 0, 0),
  @('Papua Nueva Guinea', 2, 0, 0, 2, 0, 0, 0),
  @('Bonaire, San Eustaquio y Saba', 2, 0, 0, 2, 0, 0, 0),
  @('Timor Oriental', 1, 0, 0, 1, 0, 0, 0),
  @('San Pedro y Miquelon', 1, 0, 0, 1, 0, 0, 0)
)

$data = New-Object 'object[,]' $rows.Count,8
for ($r = 0; $r -lt $rows.Count; $r++) {
  for ($c = 0; $c -lt 8; $c++) {
    $data[$r,$c] = $rows[$r][$c]
  }
}

$ws.Range('A4:H215').Value2 = $data
Write-Output "applied paises.xlsx update: $($rows.Count) rows"
